$wb = $excel.ActiveWorkbook

# Locate existing sheets
$sheetCS206 = $wb.Worksheets.Item("CS206")

# Insert a new worksheet "CS201" right after "CS206" (i.e. before "CS204")
$newSheet = $wb.Worksheets.Add($null, $sheetCS206)
$newSheet.Name = "CS201"

# Re-fetch CS204 reference now that the sheet collection changed
$sheetCS204 = $wb.Worksheets.Item("CS204")

# Fill column B (names) first, then column A (ids), to match the
# shared-string insertion order of the source workbook.
$names = @("Student7", "Student8", "Student9", "Student10")
for ($i = 0; $i -lt $names.Length; $i++) {
    $newSheet.Cells.Item($i + 1, 2).Value = $names[$i]
}

$codes = @("cs666667", "cs666677", "cs666672", "cs666673")
for ($i = 0; $i -lt $codes.Length; $i++) {
    $newSheet.Cells.Item($i + 1, 1).Value = $codes[$i]
}

# Remaining rows reuse existing students already present in the workbook
$extra = @(
    @("cs122011", "Hafiz Muhammad Haris"),
    @("cs122030", "Muhammad Touseef Khan"),
    @("cs122043", "Uroosa Shahid")
)
$row = $names.Length + 1
foreach ($pair in $extra) {
    $newSheet.Cells.Item($row, 1).Value = $pair[0]
    $newSheet.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Column B width to match source formatting (auto-sized to content, like the
# other data columns in this workbook)
$newSheet.Columns.Item(2).AutoFit()

# Update selection on CS204 (old sheet), then activate CS201 and set its selection
$sheetCS204.Activate() | Out-Null
$sheetCS204.Range("B17").Select() | Out-Null

$newSheet.Activate() | Out-Null
$newSheet.Range("E10").Select() | Out-Null
